$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Example questionnaire")

# Update ratings for rows 17 and 18 (confidence column C, rating_1to5)
$ws.Range("C17").Value = 5
$ws.Range("C18").Value = 4

# Update the active cell selection to D18
$ws.Range("D18").Select()
